$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows whose content changed (data re-shuffled across rows) ---
# Row 12
$ws.Range("F12").Value = "Csikszereda M. Ciuc"
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = "Progresul Spartac"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1.61
$ws.Range("K12").Value = "10/08/2023 22:12"
$ws.Range("L12").Value = 1.19
$ws.Range("M12").Value = "12/08/2023 09:57"
$ws.Range("N12").Value = 3.71
$ws.Range("O12").Value = "10/08/2023 22:12"
$ws.Range("P12").Value = 7.31
$ws.Range("Q12").Value = "12/08/2023 09:57"
$ws.Range("R12").Value = 4.94
$ws.Range("S12").Value = "10/08/2023 22:12"
$ws.Range("T12").Value = 13.77
$ws.Range("U12").Value = "12/08/2023 09:57"
$ws.Range("V12").Value = "https://www.betexplorer.com/football/romania/liga-2/miercurea-ciuc-progresul-spartac/v5l59JZ2/"

# Row 16
$ws.Range("F16").Value = "Selimbar"
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = "Unirea Dej"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2.06
$ws.Range("K16").Value = "10/08/2023 22:12"
$ws.Range("L16").Value = 2.2
$ws.Range("M16").Value = "12/08/2023 09:56"
$ws.Range("N16").Value = 3
$ws.Range("O16").Value = "10/08/2023 22:12"
$ws.Range("P16").Value = 2.98
$ws.Range("Q16").Value = "12/08/2023 09:56"
$ws.Range("R16").Value = 3.52
$ws.Range("S16").Value = "10/08/2023 22:12"
$ws.Range("T16").Value = 3.73
$ws.Range("U16").Value = "12/08/2023 09:56"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/romania/liga-2/selimbar-unirea-dej/OYyYABgA/"

# Row 93
$ws.Range("F93").Value = "Csikszereda M. Ciuc"
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = "Unirea Dej"
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 1.6
$ws.Range("K93").Value = "19/10/2023 21:12"
$ws.Range("L93").Value = 1.75
$ws.Range("M93").Value = "21/10/2023 09:59"
$ws.Range("N93").Value = 3.67
$ws.Range("O93").Value = "19/10/2023 21:12"
$ws.Range("P93").Value = 3.38
$ws.Range("Q93").Value = "21/10/2023 09:59"
$ws.Range("R93").Value = 4.89
$ws.Range("S93").Value = "19/10/2023 21:12"
$ws.Range("T93").Value = 5.3
$ws.Range("U93").Value = "21/10/2023 09:59"
$ws.Range("V93").Value = "https://www.betexplorer.com/football/romania/liga-2/miercurea-ciuc-unirea-dej/h6Tjc7s1/"

# Row 94
$ws.Range("F94").Value = "Progresul Spartac"
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = "Metaloglobus Bucharest"
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = 3.91
$ws.Range("K94").Value = "19/10/2023 21:12"
$ws.Range("L94").Value = 5.54
$ws.Range("M94").Value = "21/10/2023 09:58"
$ws.Range("N94").Value = 3.34
$ws.Range("O94").Value = "19/10/2023 21:12"
$ws.Range("P94").Value = 3.52
$ws.Range("Q94").Value = "21/10/2023 09:58"
$ws.Range("R94").Value = 1.83
$ws.Range("S94").Value = "19/10/2023 21:12"
$ws.Range("T94").Value = 1.68
$ws.Range("U94").Value = "21/10/2023 09:58"
$ws.Range("V94").Value = "https://www.betexplorer.com/football/romania/liga-2/progresul-spartac-metaloglobus-bucharest/nTSfdRd7/"

# Row 103
$ws.Range("F103").Value = "Unirea Dej"
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = "Progresul Spartac"
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1.61
$ws.Range("K103").Value = "26/10/2023 22:12"
$ws.Range("L103").Value = 1.62
$ws.Range("M103").Value = "28/10/2023 09:54"
$ws.Range("N103").Value = 3.6
$ws.Range("O103").Value = "26/10/2023 22:12"
$ws.Range("P103").Value = 3.64
$ws.Range("Q103").Value = "28/10/2023 09:54"
$ws.Range("R103").Value = 4.87
$ws.Range("S103").Value = "26/10/2023 22:12"
$ws.Range("T103").Value = 6.09
$ws.Range("U103").Value = "28/10/2023 09:54"
$ws.Range("V103").Value = "https://www.betexplorer.com/football/romania/liga-2/unirea-dej-progresul-spartac/0WhUAQB0/"

# Row 104
$ws.Range("F104").Value = "Tunari"
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = "Viitorul Tg. Jiu"
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 2.13
$ws.Range("K104").Value = "26/10/2023 22:12"
$ws.Range("L104").Value = 2.18
$ws.Range("M104").Value = "28/10/2023 09:59"
$ws.Range("N104").Value = 3.22
$ws.Range("O104").Value = "26/10/2023 22:12"
$ws.Range("P104").Value = 3.46
$ws.Range("Q104").Value = "28/10/2023 09:59"
$ws.Range("R104").Value = 3.11
$ws.Range("S104").Value = "26/10/2023 22:12"
$ws.Range("T104").Value = 3.24
$ws.Range("U104").Value = "28/10/2023 09:59"
$ws.Range("V104").Value = "https://www.betexplorer.com/football/romania/liga-2/tunari-viitorul-targu-jiu/CA0HDSRs/"

# Row 105
$ws.Range("F105").Value = "Selimbar"
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = "CSC Dumbravita"
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 1.75
$ws.Range("K105").Value = "26/10/2023 22:12"
$ws.Range("L105").Value = 1.79
$ws.Range("M105").Value = "28/10/2023 09:51"
$ws.Range("N105").Value = 3.4
$ws.Range("O105").Value = "26/10/2023 22:12"
$ws.Range("P105").Value = 3.56
$ws.Range("Q105").Value = "28/10/2023 09:58"
$ws.Range("R105").Value = 4.23
$ws.Range("S105").Value = "26/10/2023 22:12"
$ws.Range("T105").Value = 4.62
$ws.Range("U105").Value = "28/10/2023 09:51"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/romania/liga-2/selimbar-csc-dumbravita/dbnw94tD/"

# Row 106
$ws.Range("F106").Value = "Mioveni"
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = "Concordia"
$ws.Range("I106").Value = 1
$ws.Range("J106").Value = 3.02
$ws.Range("K106").Value = "26/10/2023 22:12"
$ws.Range("L106").Value = 3.13
$ws.Range("M106").Value = "28/10/2023 09:52"
$ws.Range("N106").Value = 3.04
$ws.Range("O106").Value = "26/10/2023 22:12"
$ws.Range("P106").Value = 2.92
$ws.Range("Q106").Value = "28/10/2023 09:09"
$ws.Range("R106").Value = 2.27
$ws.Range("S106").Value = "26/10/2023 22:12"
$ws.Range("T106").Value = 2.53
$ws.Range("U106").Value = "28/10/2023 09:52"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/romania/liga-2/mioveni-concordia/A1EkI8JQ/"

# Row 108
$ws.Range("F108").Value = "Chindia Targoviste"
$ws.Range("G108").Value = 4
$ws.Range("H108").Value = "CSM Resita"
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = 1.57
$ws.Range("K108").Value = "28/10/2023 09:25"
$ws.Range("L108").Value = 1.55
$ws.Range("M108").Value = "28/10/2023 09:59"
$ws.Range("N108").Value = 3.96
$ws.Range("O108").Value = "28/10/2023 09:25"
$ws.Range("P108").Value = 4.01
$ws.Range("Q108").Value = "28/10/2023 09:59"
$ws.Range("R108").Value = 6.01
$ws.Range("S108").Value = "28/10/2023 09:25"
$ws.Range("T108").Value = 6.19
$ws.Range("U108").Value = "28/10/2023 09:59"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/romania/liga-2/chindia-targoviste-csm-resita/fZkMCnsl/"

# Row 109
$ws.Range("F109").Value = "Alexandria"
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = "Csikszereda M. Ciuc"
$ws.Range("I109").Value = 2
$ws.Range("J109").Value = 3.7
$ws.Range("K109").Value = "26/10/2023 22:12"
$ws.Range("L109").Value = 3.93
$ws.Range("M109").Value = "28/10/2023 09:42"
$ws.Range("N109").Value = 3.19
$ws.Range("O109").Value = "26/10/2023 22:12"
$ws.Range("P109").Value = 3.24
$ws.Range("Q109").Value = "28/10/2023 09:42"
$ws.Range("R109").Value = 1.93
$ws.Range("S109").Value = "26/10/2023 22:12"
$ws.Range("T109").Value = 2.02
$ws.Range("U109").Value = "28/10/2023 09:42"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/romania/liga-2/csm-alexandria-miercurea-ciuc/vNiY9pR6/"

# Row 112
$ws.Range("F112").Value = "Ceahlaul"
$ws.Range("G112").Value = 4
$ws.Range("H112").Value = "CSM Slatina"
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2
$ws.Range("K112").Value = "04/11/2023 00:42"
$ws.Range("L112").Value = 2.04
$ws.Range("M112").Value = "04/11/2023 09:09"
$ws.Range("N112").Value = 3.17
$ws.Range("O112").Value = "04/11/2023 00:42"
$ws.Range("P112").Value = 3.17
$ws.Range("Q112").Value = "04/11/2023 09:09"
$ws.Range("R112").Value = 3.85
$ws.Range("S112").Value = "04/11/2023 00:42"
$ws.Range("T112").Value = 3.99
$ws.Range("U112").Value = "04/11/2023 09:09"
$ws.Range("V112").Value = "https://www.betexplorer.com/football/romania/liga-2/ceahlaul-csm-slatina/GK44Cwm5/"

# Row 113
$ws.Range("F113").Value = "Hunedoara"
$ws.Range("G113").Value = 4
$ws.Range("H113").Value = "Unirea Dej"
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1.43
$ws.Range("K113").Value = "02/11/2023 22:12"
$ws.Range("L113").Value = 1.66
$ws.Range("M113").Value = "04/11/2023 09:52"
$ws.Range("N113").Value = 4.09
$ws.Range("O113").Value = "02/11/2023 22:12"
$ws.Range("P113").Value = 3.67
$ws.Range("Q113").Value = "04/11/2023 09:52"
$ws.Range("R113").Value = 6.23
$ws.Range("S113").Value = "02/11/2023 22:12"
$ws.Range("T113").Value = 5.45
$ws.Range("U113").Value = "04/11/2023 09:52"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/romania/liga-2/corvinul-hunedoara-unirea-dej/p0dL8enU/"

# Row 114
$ws.Range("F114").Value = "CSC Dumbravita"
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = "Steaua Bucuresti"
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 3.3
$ws.Range("K114").Value = "02/11/2023 22:12"
$ws.Range("L114").Value = 4.09
$ws.Range("M114").Value = "04/11/2023 09:56"
$ws.Range("N114").Value = 3.34
$ws.Range("O114").Value = "02/11/2023 22:12"
$ws.Range("P114").Value = 3.61
$ws.Range("Q114").Value = "04/11/2023 09:57"
$ws.Range("R114").Value = 2
$ws.Range("S114").Value = "02/11/2023 22:12"
$ws.Range("T114").Value = 1.86
$ws.Range("U114").Value = "04/11/2023 09:57"
$ws.Range("V114").Value = "https://www.betexplorer.com/football/romania/liga-2/csc-dumbravita-csa-steaua-bucuresti/fB38Bc2B/"

# Row 115
$ws.Range("F115").Value = "CSM Resita"
$ws.Range("G115").Value = 5
$ws.Range("H115").Value = "Metaloglobus Bucharest"
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 2.38
$ws.Range("K115").Value = "02/11/2023 22:12"
$ws.Range("L115").Value = 2.59
$ws.Range("M115").Value = "04/11/2023 08:14"
$ws.Range("N115").Value = 3
$ws.Range("O115").Value = "02/11/2023 22:12"
$ws.Range("P115").Value = 3.01
$ws.Range("Q115").Value = "04/11/2023 09:52"
$ws.Range("R115").Value = 2.88
$ws.Range("S115").Value = "02/11/2023 22:12"
$ws.Range("T115").Value = 2.83
$ws.Range("U115").Value = "04/11/2023 09:54"
$ws.Range("V115").Value = "https://www.betexplorer.com/football/romania/liga-2/csm-resita-metaloglobus-bucharest/n9rm4Ffo/"

# Row 116
$ws.Range("F116").Value = "Mioveni"
$ws.Range("G116").Value = 4
$ws.Range("H116").Value = "Tunari"
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 1.55
$ws.Range("K116").Value = "02/11/2023 22:12"
$ws.Range("L116").Value = 1.69
$ws.Range("M116").Value = "04/11/2023 08:41"
$ws.Range("N116").Value = 3.8
$ws.Range("O116").Value = "02/11/2023 22:12"
$ws.Range("P116").Value = 3.65
$ws.Range("Q116").Value = "04/11/2023 08:41"
$ws.Range("R116").Value = 5.13
$ws.Range("S116").Value = "02/11/2023 22:12"
$ws.Range("T116").Value = 5.2
$ws.Range("U116").Value = "04/11/2023 08:41"
$ws.Range("V116").Value = "https://www.betexplorer.com/football/romania/liga-2/mioveni-tunari/4pDXBRDi/"

# --- Add new rows 119 and 120 (formats copied from row 118 as template) ---
# Row 119
$ws.Range("A118:V118").Copy()
$ws.Range("A119:V119").PasteSpecial(-4122)
$ws.Range("A119").Value = 118
$ws.Range("B119").Value = "romania"
$ws.Range("C119").Value = "liga-2"
$ws.Range("D119").Value = "2023-2024"
$ws.Range("E119").Value = 45235.41666666666
$ws.Range("F119").Value = "Progresul Spartac"
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = "Alexandria"
$ws.Range("I119").Value = 3
$ws.Range("J119").Value = 3.4
$ws.Range("K119").Value = "04/11/2023 22:41"
$ws.Range("L119").Value = 3.4
$ws.Range("M119").Value = "04/11/2023 22:41"
$ws.Range("N119").Value = 3.36
$ws.Range("O119").Value = "05/11/2023 08:01"
$ws.Range("P119").Value = 3.36
$ws.Range("Q119").Value = "05/11/2023 08:01"
$ws.Range("R119").Value = 2.12
$ws.Range("S119").Value = "04/11/2023 22:41"
$ws.Range("T119").Value = 2.12
$ws.Range("U119").Value = "04/11/2023 22:41"
$ws.Range("V119").Value = "https://www.betexplorer.com/football/romania/liga-2/progresul-spartac-csm-alexandria/jwcH9yXN/"

# Row 120
$ws.Range("A118:V118").Copy()
$ws.Range("A120:V120").PasteSpecial(-4122)
$ws.Range("A120").Value = 119
$ws.Range("B120").Value = "romania"
$ws.Range("C120").Value = "liga-2"
$ws.Range("D120").Value = "2023-2024"
$ws.Range("E120").Value = 45235.47916666666
$ws.Range("F120").Value = "Csikszereda M. Ciuc"
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = "Selimbar"
$ws.Range("I120").Value = 2
$ws.Range("J120").Value = 1.97
$ws.Range("K120").Value = "03/11/2023 23:42"
$ws.Range("L120").Value = 1.99
$ws.Range("M120").Value = "05/11/2023 11:21"
$ws.Range("N120").Value = 3.08
$ws.Range("O120").Value = "03/11/2023 23:42"
$ws.Range("P120").Value = 3.21
$ws.Range("Q120").Value = "05/11/2023 11:21"
$ws.Range("R120").Value = 3.73
$ws.Range("S120").Value = "03/11/2023 23:42"
$ws.Range("T120").Value = 4.1
$ws.Range("U120").Value = "05/11/2023 11:21"
$ws.Range("V120").Value = "https://www.betexplorer.com/football/romania/liga-2/miercurea-ciuc-selimbar/xzgDAHHH/"

$excel.CutCopyMode = 0